# Add a new worksheet "resumen hogares" with a summary table of households,
# placed after the existing "Visistas" sheet (i.e. at the end).
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "resumen hogares"

# Headers
$ws.Range("A1").Value = "Id_hogar"
$ws.Range("B1").Value = "Edad_promedio"
$ws.Range("C1").Value = "Genero"
$ws.Range("D1").Value = "Nivel_educacional"
$ws.Range("E1").Value = "Estado_civil"
$ws.Range("F1").Value = "Nacionalidad"

# Data rows
$ws.Range("A2").Value = 33
$ws.Range("B2").Value = 46
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2

$ws.Range("A3").Value = 35
$ws.Range("B3").Value = 29
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 4

$ws.Range("A4").Value = 36
$ws.Range("B4").Value = 23
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 2
